$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2 values: B2 -> "Gold", D2 -> 0.337 (A2 and C2 stay 0)
$ws.Range("B2").Value = "Gold"
$ws.Range("D2").Value = 0.337

# Delete rows 3 through 7 (no longer present in the target sheet)
$ws.Range("A3:D7").EntireRow.Delete()
